# Update FFXIV market-price/profit figures (columns H-N) across all leve-crafting
# sheets, refreshed from the latest scheduled market-data run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43: Growing Is Knowing / Growth Formula Gamma
$ws.Range("H43").Value = 5324.1665
$ws.Range("I43").Value = 5237.5
$ws.Range("J43").Value = 5497.5
$ws.Range("K43").Value = 5237.5
$ws.Range("L43").Value = 5497.5
$ws.Range("M43").Value = -5168.5
$ws.Range("N43").Value = -5635.5

# Row 62: The Mustache Suits Him / Enchanted Mythrite Ink
$ws.Range("H62").Value = 3942.6
$ws.Range("I62").Value = 1453.9
$ws.Range("K62").Value = 1453.9
$ws.Range("M62").Value = -829.9000000000001

# Row 65: Forgery of Convenience (L) / Enchanted Mythrite Ink
$ws.Range("H65").Value = 3942.6
$ws.Range("I65").Value = 1453.9
$ws.Range("K65").Value = 7269.5
$ws.Range("M65").Value = -4149.5

# Row 107: Another Man's Ink / Enchanted Truegold Ink
$ws.Range("H107").Value = 962.2353000000001
$ws.Range("I107").Value = 555.9167
$ws.Range("K107").Value = 555.9167
$ws.Range("M107").Value = 1364.0833

# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 8162.4375
$ws.Range("I116").Value = 7447.8335
$ws.Range("K116").Value = 7447.8335
$ws.Range("M116").Value = -4005.8335

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 4077
$ws.Range("J2").Value = 3333.3333
$ws.Range("L2").Value = 3333.3333
$ws.Range("N2").Value = -3559.3333

# Row 32: Ingot We Trust / Steel Ingot
$ws.Range("H32").Value = 4666.087
$ws.Range("I32").Value = 4730.2925
$ws.Range("K32").Value = 4730.2925
$ws.Range("M32").Value = -4443.2925

# Row 61: Dealing with the Tough Stuff / Cobalt Ingot
$ws.Range("H61").Value = 7243397
$ws.Range("I61").Value = 7144946.5
$ws.Range("K61").Value = 7144946.5
$ws.Range("M61").Value = -7144734.5

# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 963338.6
$ws.Range("I74").Value = 1042866.9
$ws.Range("J74").Value = 8999.5
$ws.Range("K74").Value = 1042866.9
$ws.Range("L74").Value = 8999.5
$ws.Range("M74").Value = -1041992.9
$ws.Range("N74").Value = -10747.5

# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 963338.6
$ws.Range("I77").Value = 1042866.9
$ws.Range("J77").Value = 8999.5
$ws.Range("K77").Value = 5214334.5
$ws.Range("L77").Value = 44997.5
$ws.Range("M77").Value = -5209966.5
$ws.Range("N77").Value = -53733.5

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 4840
$ws.Range("I102").Value = 3741
$ws.Range("K102").Value = 3741
$ws.Range("M102").Value = -2119

# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 8233.444
$ws.Range("I110").Value = 7728.7144
$ws.Range("K110").Value = 7728.7144
$ws.Range("M110").Value = -5683.7144

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 4077
$ws.Range("J116").Value = 3333.3333
$ws.Range("L116").Value = 3333.3333
$ws.Range("N116").Value = -7921.3333

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 3723.9412
$ws.Range("I122").Value = 3531.5
$ws.Range("K122").Value = 10594.5
$ws.Range("M122").Value = -8144.5

# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 1495355.9
$ws.Range("I132").Value = 2559.3447
$ws.Range("K132").Value = 7678.034100000001
$ws.Range("M132").Value = -5148.034100000001

# Row 136: Metal with Mettle / Cobalt Tungsten Ingot
$ws.Range("H136").Value = 7243397
$ws.Range("I136").Value = 7144946.5
$ws.Range("K136").Value = 21434839.5
$ws.Range("M136").Value = -21432289.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 4077
$ws.Range("J3").Value = 3333.3333
$ws.Range("L3").Value = 3333.3333
$ws.Range("N3").Value = -3561.3333

# Row 20: Smelt and Dealt / Iron Ingot
$ws.Range("H20").Value = 7788.4287
$ws.Range("J20").Value = 2974
$ws.Range("L20").Value = 2974
$ws.Range("N20").Value = -3468

# Row 80: Unbreaker / Titanium Ingot
$ws.Range("H80").Value = 912.61536
$ws.Range("I80").Value = 598.6667
$ws.Range("K80").Value = 598.6667
$ws.Range("M80").Value = 399.3333

# Row 83: Attack on Titanium (L) / Titanium Ingot
$ws.Range("H83").Value = 912.61536
$ws.Range("I83").Value = 598.6667
$ws.Range("K83").Value = 2993.3335
$ws.Range("M83").Value = 1998.6665

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 4168785.8
$ws.Range("I134").Value = 2041
$ws.Range("K134").Value = 6123
$ws.Range("M134").Value = -3588

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof / Ash Lumber
$ws.Range("H16").Value = 4768022
$ws.Range("I16").Value = 7147910
$ws.Range("K16").Value = 7147910
$ws.Range("M16").Value = -7147623

# Row 99: O Pine / Pine Lumber
$ws.Range("H99").Value = 22809.818
$ws.Range("I99").Value = 12499.167
$ws.Range("J99").Value = 35182.6
$ws.Range("K99").Value = 12499.167
$ws.Range("L99").Value = 35182.6
$ws.Range("M99").Value = -11001.167
$ws.Range("N99").Value = -38178.6

# Row 113: Patient Patients / White Ash Lumber
$ws.Range("H113").Value = 4768022
$ws.Range("I113").Value = 7147910
$ws.Range("K113").Value = 7147910
$ws.Range("M113").Value = -7145740

# Row 126: A Better Conductor / Red Pine Lumber
$ws.Range("H126").Value = 22809.818
$ws.Range("I126").Value = 12499.167
$ws.Range("J126").Value = 35182.6
$ws.Range("K126").Value = 37497.501
$ws.Range("L126").Value = 105547.8
$ws.Range("M126").Value = -35027.501
$ws.Range("N126").Value = -110487.8

$ws = $wb.Worksheets.Item("CUL")
# Row 107: Slippery Service / Frantoio Oil
$ws.Range("H107").Value = 5358042
$ws.Range("I107").Value = 3815.6
$ws.Range("J107").Value = 7588969.5
$ws.Range("K107").Value = 11446.8
$ws.Range("L107").Value = 22766908.5
$ws.Range("M107").Value = -9526.799999999999
$ws.Range("N107").Value = -22770748.5

# Row 130: Blast from the Pasta / The Noodles of Elpis
$ws.Range("H130").Value = 11256.889
$ws.Range("I130").Value = 3326.6667
$ws.Range("K130").Value = 9980.000100000001
$ws.Range("M130").Value = -4960.000100000001

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 4104.1
$ws.Range("J131").Value = 7088.5713
$ws.Range("L131").Value = 21265.7139
$ws.Range("N131").Value = -31345.7139

# Row 138: Bring Me Your Tacos / Tacos Al Pastor
$ws.Range("H138").Value = 10498.786
$ws.Range("I138").Value = 9360.833000000001
$ws.Range("K138").Value = 28082.499
$ws.Range("M138").Value = -22942.499

# Row 139: Najoothie / Wild Banana Blend
$ws.Range("H139").Value = 4520.6665
$ws.Range("I139").Value = 3080.05
$ws.Range("K139").Value = 9240.150000000001
$ws.Range("M139").Value = -4100.150000000001

$ws = $wb.Worksheets.Item("GSM")
# Row 107: Whetstones for the Workers / Hard Mudstone Whetstone
$ws.Range("H107").Value = 2851
$ws.Range("I107").Value = 2851
$ws.Range("K107").Value = 2851
$ws.Range("M107").Value = -931

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 4338.077
$ws.Range("I122").Value = 4650.25
$ws.Range("J122").Value = 3838.6
$ws.Range("K122").Value = 13950.75
$ws.Range("L122").Value = 11515.8
$ws.Range("M122").Value = -11500.75
$ws.Range("N122").Value = -16415.8

# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 2074.125
$ws.Range("I126").Value = 1712.4
$ws.Range("K126").Value = 5137.200000000001
$ws.Range("M126").Value = -2667.200000000001

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 4765161.5
$ws.Range("J132").Value = 12503607
$ws.Range("L132").Value = 37510821
$ws.Range("N132").Value = -37515881

$ws = $wb.Worksheets.Item("LTW")
# Row 22: Skin off Their Backs / Aldgoat Leather
$ws.Range("H22").Value = 11506.4
$ws.Range("J22").Value = 7749.25
$ws.Range("L22").Value = 7749.25
$ws.Range("N22").Value = -8339.25

# Row 27: Fire and Hide / Aldgoat Leather
$ws.Range("H27").Value = 11506.4
$ws.Range("J27").Value = 7749.25
$ws.Range("L27").Value = 7749.25
$ws.Range("N27").Value = -7963.25

# Row 46: Supply Side Logic / Boar Leather
$ws.Range("H46").Value = 8599
$ws.Range("J46").Value = 7500
$ws.Range("L46").Value = 7500
$ws.Range("N46").Value = -7876

# Row 55: It's Not a Job, It's a Calling / Peiste Leather
$ws.Range("H55").Value = 1450.9286
$ws.Range("J55").Value = 1246.75
$ws.Range("L55").Value = 1246.75
$ws.Range("N55").Value = -1592.75

# Row 61: Spelling Me Softly / Raptor Leather
$ws.Range("H61").Value = 10833
$ws.Range("J61").Value = 17668.334
$ws.Range("L61").Value = 17668.334
$ws.Range("N61").Value = -18072.334

# Row 113: Peace in Rest / Atrociraptor Leather
$ws.Range("H113").Value = 10833
$ws.Range("J113").Value = 17668.334
$ws.Range("L113").Value = 17668.334
$ws.Range("N113").Value = -22008.334

# Row 122: Hell on Leather / Gaja Leather
$ws.Range("H122").Value = 3375.14
$ws.Range("I122").Value = 3233.1702
$ws.Range("J122").Value = 5599.3335
$ws.Range("K122").Value = 9699.5106
$ws.Range("L122").Value = 16798.0005
$ws.Range("M122").Value = -7249.5106
$ws.Range("N122").Value = -21698.0005

$ws = $wb.Worksheets.Item("WVR")
# Row 135: In Line with Linen / Mountain Linen Cloak of Casting
$ws.Range("H135").Value = 101172.5
$ws.Range("J135").Value = 101172.5
$ws.Range("L135").Value = 101172.5
$ws.Range("N135").Value = -111312.5

# Row 136: Weaving the Envelope / Sarcenet Cloth
$ws.Range("H136").Value = 234638.77
$ws.Range("I136").Value = 1933.7949
$ws.Range("K136").Value = 5801.384700000001
$ws.Range("M136").Value = -3251.384700000001
